$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from G4 onto G5 before writing values there,
# so the new cell reuses the existing style (s="1") instead of creating a new one.
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5").Value = 9975.93
$ws.Range("B5").Value = 10017
$ws.Range("C5").Value = 80.11
$ws.Range("D5").Value = 79.78
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = -0.41
$ws.Range("G5").Value = 42609.505231481482
$ws.Range("H5").Value = $false
